$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1616.25
$ws.Cells.Item(17, 10).Value = 1616.25
$ws.Cells.Item(17, 12).Value = 4848.75
$ws.Cells.Item(17, 14).Value = -5184.75
$ws.Cells.Item(40, 8).Value = 2950
$ws.Cells.Item(40, 9).Value = 2920
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 2920
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = -2745
$ws.Cells.Item(40, 14).Value = -3350
$ws.Cells.Item(100, 8).Value = 640.6
$ws.Cells.Item(100, 9).Value = 677.25
$ws.Cells.Item(100, 10).Value = 494
$ws.Cells.Item(100, 11).Value = 677.25
$ws.Cells.Item(100, 12).Value = 494
$ws.Cells.Item(100, 13).Value = -136.25
$ws.Cells.Item(100, 14).Value = -1576
$ws.Cells.Item(132, 8).Value = 385806.78
$ws.Cells.Item(132, 9).Value = 1214.95
$ws.Cells.Item(132, 10).Value = 1667779.5
$ws.Cells.Item(132, 11).Value = 3644.85
$ws.Cells.Item(132, 12).Value = 5003338.5
$ws.Cells.Item(132, 13).Value = -1114.85
$ws.Cells.Item(132, 14).Value = -5008398.5
$ws.Cells.Item(137, 8).Value = 3024.8948
$ws.Cells.Item(137, 9).Value = 1731.75
$ws.Cells.Item(137, 10).Value = 9921.666999999999
$ws.Cells.Item(137, 11).Value = 5195.25
$ws.Cells.Item(137, 12).Value = 29765.001
$ws.Cells.Item(137, 13).Value = -2645.25
$ws.Cells.Item(137, 14).Value = -34865.001
$ws.Cells.Item(138, 8).Value = 4452.2095
$ws.Cells.Item(138, 9).Value = 2544.6875
$ws.Cells.Item(138, 10).Value = 5582.593
$ws.Cells.Item(138, 11).Value = 7634.0625
$ws.Cells.Item(138, 12).Value = 16747.779
$ws.Cells.Item(138, 13).Value = -2494.0625
$ws.Cells.Item(138, 14).Value = -27027.779
$ws.Cells.Item(141, 8).Value = 4410.4
$ws.Cells.Item(141, 9).Value = 3783.6155
$ws.Cells.Item(141, 10).Value = 6221.1113
$ws.Cells.Item(141, 11).Value = 11350.8465
$ws.Cells.Item(141, 12).Value = 18663.3339
$ws.Cells.Item(141, 13).Value = -6170.8465
$ws.Cells.Item(141, 14).Value = -29023.3339

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4154.727
$ws.Cells.Item(32, 9).Value = 2006.2394
$ws.Cells.Item(32, 11).Value = 2006.2394
$ws.Cells.Item(32, 13).Value = -1719.2394
$ws.Cells.Item(61, 8).Value = 2401.2
$ws.Cells.Item(61, 9).Value = 1864.0454
$ws.Cells.Item(61, 11).Value = 1864.0454
$ws.Cells.Item(61, 13).Value = -1652.0454
$ws.Cells.Item(102, 8).Value = 1948.1724
$ws.Cells.Item(102, 9).Value = 1059.92
$ws.Cells.Item(102, 10).Value = 7499.75
$ws.Cells.Item(102, 11).Value = 1059.92
$ws.Cells.Item(102, 12).Value = 7499.75
$ws.Cells.Item(102, 13).Value = 562.0799999999999
$ws.Cells.Item(102, 14).Value = -10743.75
$ws.Cells.Item(132, 8).Value = 2168.1538
$ws.Cells.Item(132, 9).Value = 2158.1304
$ws.Cells.Item(132, 10).Value = 2245
$ws.Cells.Item(132, 11).Value = 6474.3912
$ws.Cells.Item(132, 12).Value = 6735
$ws.Cells.Item(132, 13).Value = -3944.3912
$ws.Cells.Item(132, 14).Value = -11795
$ws.Cells.Item(136, 8).Value = 2401.2
$ws.Cells.Item(136, 9).Value = 1864.0454
$ws.Cells.Item(136, 11).Value = 5592.1362
$ws.Cells.Item(136, 13).Value = -3042.1362

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 626.6429000000001
$ws.Cells.Item(80, 9).Value = 67.14286
$ws.Cells.Item(80, 11).Value = 67.14286
$ws.Cells.Item(80, 13).Value = 930.85714
$ws.Cells.Item(83, 8).Value = 626.6429000000001
$ws.Cells.Item(83, 9).Value = 67.14286
$ws.Cells.Item(83, 11).Value = 335.7143
$ws.Cells.Item(83, 13).Value = 4656.2857
$ws.Cells.Item(134, 8).Value = 2636.7104
$ws.Cells.Item(134, 9).Value = 2592.4375
$ws.Cells.Item(134, 10).Value = 2872.8333
$ws.Cells.Item(134, 11).Value = 7777.3125
$ws.Cells.Item(134, 12).Value = 8618.499899999999
$ws.Cells.Item(134, 13).Value = -5242.3125
$ws.Cells.Item(134, 14).Value = -13688.4999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7173.213
$ws.Cells.Item(31, 9).Value = 5097.1055
$ws.Cells.Item(31, 10).Value = 8582
$ws.Cells.Item(31, 11).Value = 5097.1055
$ws.Cells.Item(31, 12).Value = 8582
$ws.Cells.Item(31, 13).Value = -4802.1055
$ws.Cells.Item(31, 14).Value = -9172
$ws.Cells.Item(34, 8).Value = 7173.213
$ws.Cells.Item(34, 9).Value = 5097.1055
$ws.Cells.Item(34, 10).Value = 8582
$ws.Cells.Item(34, 11).Value = 5097.1055
$ws.Cells.Item(34, 12).Value = 8582
$ws.Cells.Item(34, 13).Value = -4895.1055
$ws.Cells.Item(34, 14).Value = -8986
$ws.Cells.Item(58, 8).Value = 2274.2942
$ws.Cells.Item(58, 9).Value = 2304.0625
$ws.Cells.Item(58, 11).Value = 2304.0625
$ws.Cells.Item(58, 13).Value = -2101.0625
$ws.Cells.Item(111, 8).Value = 75000
$ws.Cells.Item(111, 10).Value = 75000
$ws.Cells.Item(111, 12).Value = 75000
$ws.Cells.Item(111, 14).Value = -83180
$ws.Cells.Item(122, 8).Value = 5029.5
$ws.Cells.Item(122, 9).Value = 4932.778
$ws.Cells.Item(122, 11).Value = 14798.334
$ws.Cells.Item(122, 13).Value = -12348.334
$ws.Cells.Item(132, 8).Value = 1452.75
$ws.Cells.Item(132, 9).Value = 1264.5454
$ws.Cells.Item(132, 10).Value = 2142.8333
$ws.Cells.Item(132, 11).Value = 3793.6362
$ws.Cells.Item(132, 12).Value = 6428.499899999999
$ws.Cells.Item(132, 13).Value = -1263.6362
$ws.Cells.Item(132, 14).Value = -11488.4999
$ws.Cells.Item(134, 8).Value = 3199.32
$ws.Cells.Item(134, 9).Value = 2109.9575
$ws.Cells.Item(134, 10).Value = 20266
$ws.Cells.Item(134, 11).Value = 6329.872499999999
$ws.Cells.Item(134, 12).Value = 60798
$ws.Cells.Item(134, 13).Value = -3794.872499999999
$ws.Cells.Item(134, 14).Value = -65868
$ws.Cells.Item(136, 8).Value = 2274.2942
$ws.Cells.Item(136, 9).Value = 2304.0625
$ws.Cells.Item(136, 11).Value = 6912.1875
$ws.Cells.Item(136, 13).Value = -4362.1875

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 277815.34
$ws.Cells.Item(2, 9).Value = 454579.72
$ws.Cells.Item(2, 10).Value = 42.714287
$ws.Cells.Item(2, 11).Value = 2727478.32
$ws.Cells.Item(2, 12).Value = 256.285722
$ws.Cells.Item(2, 13).Value = -2727365.32
$ws.Cells.Item(2, 14).Value = -482.285722
$ws.Cells.Item(29, 8).Value = 242.33333
$ws.Cells.Item(29, 9).Value = 109.8
$ws.Cells.Item(29, 10).Value = 408
$ws.Cells.Item(29, 11).Value = 329.4
$ws.Cells.Item(29, 12).Value = 1224
$ws.Cells.Item(29, 13).Value = -52.39999999999998
$ws.Cells.Item(29, 14).Value = -1778
$ws.Cells.Item(64, 8).Value = 770.25
$ws.Cells.Item(64, 9).Value = 728.6667
$ws.Cells.Item(64, 10).Value = 895
$ws.Cells.Item(64, 11).Value = 2186.0001
$ws.Cells.Item(64, 12).Value = 2685
$ws.Cells.Item(64, 13).Value = -1916.0001
$ws.Cells.Item(64, 14).Value = -3225
$ws.Cells.Item(67, 8).Value = 770.25
$ws.Cells.Item(67, 9).Value = 728.6667
$ws.Cells.Item(67, 10).Value = 895
$ws.Cells.Item(67, 11).Value = 2186.0001
$ws.Cells.Item(67, 12).Value = 2685
$ws.Cells.Item(67, 13).Value = -1250.0001
$ws.Cells.Item(67, 14).Value = -4557
$ws.Cells.Item(131, 8).Value = 4223946
$ws.Cells.Item(131, 10).Value = 4631465.5
$ws.Cells.Item(131, 12).Value = 13894396.5
$ws.Cells.Item(131, 14).Value = -13904476.5
$ws.Cells.Item(141, 8).Value = 1311.6666
$ws.Cells.Item(141, 9).Value = 1311.6666
$ws.Cells.Item(141, 11).Value = 3934.9998
$ws.Cells.Item(141, 13).Value = 1245.0002

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2436.7778
$ws.Cells.Item(132, 9).Value = 2652.6
$ws.Cells.Item(132, 11).Value = 7957.799999999999
$ws.Cells.Item(132, 13).Value = -5427.799999999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4106.875
$ws.Cells.Item(7, 9).Value = 4171.4165
$ws.Cells.Item(7, 10).Value = 3913.25
$ws.Cells.Item(7, 11).Value = 4171.4165
$ws.Cells.Item(7, 12).Value = 3913.25
$ws.Cells.Item(7, 13).Value = -4059.4165
$ws.Cells.Item(7, 14).Value = -4137.25
$ws.Cells.Item(55, 8).Value = 425.58334
$ws.Cells.Item(55, 9).Value = 364.75
$ws.Cells.Item(55, 10).Value = 547.25
$ws.Cells.Item(55, 11).Value = 364.75
$ws.Cells.Item(55, 12).Value = 547.25
$ws.Cells.Item(55, 13).Value = -191.75
$ws.Cells.Item(55, 14).Value = -893.25
$ws.Cells.Item(126, 8).Value = 4106.875
$ws.Cells.Item(126, 9).Value = 4171.4165
$ws.Cells.Item(126, 10).Value = 3913.25
$ws.Cells.Item(126, 11).Value = 12514.2495
$ws.Cells.Item(126, 12).Value = 11739.75
$ws.Cells.Item(126, 13).Value = -10044.2495
$ws.Cells.Item(126, 14).Value = -16679.75
$ws.Cells.Item(132, 8).Value = 2978.5334
$ws.Cells.Item(132, 9).Value = 2667.8
$ws.Cells.Item(132, 11).Value = 8003.400000000001
$ws.Cells.Item(132, 13).Value = -5473.400000000001
$ws.Cells.Item(136, 8).Value = 1150.9166
$ws.Cells.Item(136, 9).Value = 973
$ws.Cells.Item(136, 10).Value = 1400
$ws.Cells.Item(136, 11).Value = 2919
$ws.Cells.Item(136, 12).Value = 4200
$ws.Cells.Item(136, 13).Value = -369
$ws.Cells.Item(136, 14).Value = -9300

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 476.375
$ws.Cells.Item(100, 9).Value = 468
$ws.Cells.Item(100, 10).Value = 501.5
$ws.Cells.Item(100, 11).Value = 936
$ws.Cells.Item(100, 12).Value = 1003
$ws.Cells.Item(100, 13).Value = -395
$ws.Cells.Item(100, 14).Value = -2085
$ws.Cells.Item(126, 8).Value = 3055.2354
$ws.Cells.Item(126, 9).Value = 2562.6667
$ws.Cells.Item(126, 11).Value = 7688.000100000001
$ws.Cells.Item(126, 13).Value = -5218.000100000001
$ws.Cells.Item(132, 8).Value = 3682.7693
$ws.Cells.Item(132, 9).Value = 3221.9656
$ws.Cells.Item(132, 11).Value = 9665.8968
$ws.Cells.Item(132, 13).Value = -7135.8968
$ws.Cells.Item(136, 8).Value = 774.4286
$ws.Cells.Item(136, 9).Value = 779.7692
$ws.Cells.Item(136, 10).Value = 705
$ws.Cells.Item(136, 11).Value = 2339.3076
$ws.Cells.Item(136, 12).Value = 2115
$ws.Cells.Item(136, 13).Value = 210.6923999999999
$ws.Cells.Item(136, 14).Value = -7215
